$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "530÷6=88, 2" "801÷4=200, 1"
Replace-Text "173÷2=86, 1" "195÷3=65, 0"
Replace-Text "570÷8=71, 2" "770÷3=256, 2"
Replace-Text "753÷4=188, 1" "156÷3=52, 0"
Replace-Text "330÷3=110, 0" "240÷6=40, 0"
Replace-Text "956÷4=239, 0" "186÷7=26, 4"
Replace-Text "504÷7=72, 0" "796÷7=113, 5"
Replace-Text "841÷4=210, 1" "812÷8=101, 4"
Replace-Text "300÷5=60, 0" "228÷4=57, 0"
Replace-Text "632÷6=105, 2" "757÷5=151, 2"
Replace-Text "327÷8=40, 7" "360÷4=90, 0"
Replace-Text "668÷9=74, 2" "177÷6=29, 3"
Replace-Text "207÷2=103, 1" "872÷4=218, 0"
Replace-Text "515÷8=64, 3" "735÷9=81, 6"
Replace-Text "319÷3=106, 1" "650÷9=72, 2"
Replace-Text "348÷5=69, 3" "835÷5=167, 0"
Replace-Text "882÷3=294, 0" "403÷8=50, 3"
Replace-Text "445÷7=63, 4" "959÷2=479, 1"
Replace-Text "307÷6=51, 1" "612÷4=153, 0"
Replace-Text "857÷4=214, 1" "709÷2=354, 1"
Replace-Text "569÷5=113, 4" "541÷8=67, 5"
Replace-Text "581÷8=72, 5" "661÷8=82, 5"
Replace-Text "691÷4=172, 3" "319÷5=63, 4"
Replace-Text "261÷7=37, 2" "782÷2=391, 0"
Replace-Text "675÷8=84, 3" "139÷5=27, 4"

Write-Output "Done"
